{"js": "// Replace each old answer/date string with its corresponding new value.\n// Every `oldText` below is unique in the document, so a single search\n// hit is expected for each entry (matches the commit's 101 text edits:\n// 1 date line + 100 arithmetic-table cells).\nconst replacements = [\n  [\"2025-06-07 Saturday\", \"2025-06-08 Sunday\"],\n  [\"88-70=18\", \"2+67=69\"],\n  [\"5+81=86\", \"44+32=76\"],\n  [\"9+10=19\", \"77+9=86\"],\n  [\"88+8=96\", \"51-22=29\"],\n  [\"67-65=2\", \"13+22=35\"],\n  [\"17+75=92\", \"48+50=98\"],\n  [\"77-47=30\", \"60-7=53\"],\n  [\"3+88=91\", \"8+18=26\"],\n  [\"95-66=29\", \"3+22=25\"],\n  [\"23+4=27\", \"96-28=68\"],\n  [\"8+80=88\", \"85-9=76\"],\n  [\"38+18=56\", \"79+11=90\"],\n  [\"54-27=27\", \"49+17=66\"],\n  [\"92-30=62\", \"77-18=59\"],\n  [\"76-42=34\", \"33-16=17\"],\n  [\"2+78=80\", \"39+56=95\"],\n  [\"31+68=99\", \"36-22=14\"],\n  [\"53-11=42\", \"75-51=24\"],\n  [\"11+54=65\", \"87-64=23\"],\n  [\"1+80=81\", \"15+81=96\"],\n  [\"69+0=69\", \"85-54=31\"],\n  [\"76-37=39\", \"84+5=89\"],\n  [\"12+68=80\", \"73+7=80\"],\n  [\"41+31=72\", \"64-35=29\"],\n  [\"62-21=41\", \"19+80=99\"],\n  [\"39+3=42\", \"30-15=15\"],\n  [\"21+67=88\", \"28-7=21\"],\n  [\"31+37=68\", \"33-21=12\"],\n  [\"12+4=16\", \"29+21=50\"],\n  [\"47+32=79\", \"78-68=10\"],\n  [\"38+8=46\", \"18+30=48\"],\n  [\"25-11=14\", \"42-32=10\"],\n  [\"47+41=88\", \"54-28=26\"],\n  [\"91-41=50\", \"82+10=92\"],\n  [\"96-8=88\", \"62+30=92\"],\n  [\"33+54=87\", \"89+3=92\"],\n  [\"82-58=24\", \"45+45=90\"],\n  [\"49+11=60\", \"32+35=67\"],\n  [\"10+86=96\", \"72-8=64\"],\n  [\"71+19=90\", \"40+14=54\"],\n  [\"45+32=77\", \"50-43=7\"],\n  [\"10+53=63\", \"82-32=50\"],\n  [\"81-16=65\", \"56-8=48\"],\n  [\"98-41=57\", \"92-40=52\"],\n  [\"24+1=25\", \"15+52=67\"],\n  [\"6+85=91\", \"19+11=30\"],\n  [\"68-11=57\", \"8+34=42\"],\n  [\"48-19=29\", \"75-63=12\"],\n  [\"89-46=43\", \"52-8=44\"],\n  [\"66-32=34\", \"27+32=59\"],\n  [\"69-1=68\", \"19+66=85\"],\n  [\"12+87=99\", \"84-73=11\"],\n  [\"60+27=87\", \"90-53=37\"],\n  [\"28+71=99\", \"43-9=34\"],\n  [\"87-74=13\", \"62+8=70\"],\n  [\"61-6=55\", \"13+83=96\"],\n  [\"61+3=64\", \"16+4=20\"],\n  [\"81-45=36\", \"74+12=86\"],\n  [\"49-35=14\", \"93-73=20\"],\n  [\"74+5=79\", \"84-82=2\"],\n  [\"26-7=19\", \"27+28=55\"],\n  [\"39+46=85\", \"59+24=83\"],\n  [\"3+44=47\", \"12+1=13\"],\n  [\"56-0=56\", \"81+2=83\"],\n  [\"1+13=14\", \"53+35=88\"],\n  [\"67-60=7\", \"45+30=75\"],\n  [\"80-4=76\", \"10+67=77\"],\n  [\"63-41=22\", \"41-12=29\"],\n  [\"70-32=38\", \"4+15=19\"],\n  [\"82+11=93\", \"62-14=48\"],\n  [\"95-38=57\", \"24+10=34\"],\n  [\"77+10=87\", \"53+25=78\"],\n  [\"57-0=57\", \"61+22=83\"],\n  [\"6+6=12\", \"30-10=20\"],\n  [\"4+76=80\", \"36+32=68\"],\n  [\"1+87=88\", \"37-11=26\"],\n  [\"85+11=96\", \"56-36=20\"],\n  [\"23+5=28\", \"73+13=86\"],\n  [\"87-75=12\", \"12+54=66\"],\n  [\"55+14=69\", \"90-34=56\"],\n  [\"37+17=54\", \"91-75=16\"],\n  [\"96-5=91\", \"70+18=88\"],\n  [\"21+45=66\", \"17+4=21\"],\n  [\"23+28=51\", \"74-48=26\"],\n  [\"95-75=20\", \"22+50=72\"],\n  [\"10+44=54\", \"35+29=64\"],\n  [\"88-47=41\", \"54-35=19\"],\n  [\"17+32=49\", \"64-33=31\"],\n  [\"49-33=16\", \"68-50=18\"],\n  [\"42-21=21\", \"0+11=11\"],\n  [\"68-41=27\", \"74-10=64\"],\n  [\"57-54=3\", \"77-70=7\"],\n  [\"87-34=53\", \"88-72=16\"],\n  [\"7+60=67\", \"95-14=81\"],\n  [\"36+12=48\", \"24+14=38\"],\n  [\"87+1=88\", \"30+57=87\"],\n  [\"48+35=83\", \"86-67=19\"],\n  [\"47+49=96\", \"90-0=90\"],\n  [\"96-92=4\", \"9+42=51\"],\n  [\"58+37=95\", \"95-11=84\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replace each old answer/date string with its corresponding new value.\n# Every old string below is unique in the document, so Find.Execute with\n# Wrap=wdFindStop (1) and Replace=wdReplaceOne (2), run once per pair\n# against the whole document range, performs exactly one substitution\n# each (matches the commit's 101 text edits: 1 date line + 100\n# arithmetic-table cells).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-06-07 Saturday\", \"2025-06-08 Sunday\"),\n    @(\"88-70=18\", \"2+67=69\"),\n    @(\"5+81=86\", \"44+32=76\"),\n    @(\"9+10=19\", \"77+9=86\"),\n    @(\"88+8=96\", \"51-22=29\"),\n    @(\"67-65=2\", \"13+22=35\"),\n    @(\"17+75=92\", \"48+50=98\"),\n    @(\"77-47=30\", \"60-7=53\"),\n    @(\"3+88=91\", \"8+18=26\"),\n    @(\"95-66=29\", \"3+22=25\"),\n    @(\"23+4=27\", \"96-28=68\"),\n    @(\"8+80=88\", \"85-9=76\"),\n    @(\"38+18=56\", \"79+11=90\"),\n    @(\"54-27=27\", \"49+17=66\"),\n    @(\"92-30=62\", \"77-18=59\"),\n    @(\"76-42=34\", \"33-16=17\"),\n    @(\"2+78=80\", \"39+56=95\"),\n    @(\"31+68=99\", \"36-22=14\"),\n    @(\"53-11=42\", \"75-51=24\"),\n    @(\"11+54=65\", \"87-64=23\"),\n    @(\"1+80=81\", \"15+81=96\"),\n    @(\"69+0=69\", \"85-54=31\"),\n    @(\"76-37=39\", \"84+5=89\"),\n    @(\"12+68=80\", \"73+7=80\"),\n    @(\"41+31=72\", \"64-35=29\"),\n    @(\"62-21=41\", \"19+80=99\"),\n    @(\"39+3=42\", \"30-15=15\"),\n    @(\"21+67=88\", \"28-7=21\"),\n    @(\"31+37=68\", \"33-21=12\"),\n    @(\"12+4=16\", \"29+21=50\"),\n    @(\"47+32=79\", \"78-68=10\"),\n    @(\"38+8=46\", \"18+30=48\"),\n    @(\"25-11=14\", \"42-32=10\"),\n    @(\"47+41=88\", \"54-28=26\"),\n    @(\"91-41=50\", \"82+10=92\"),\n    @(\"96-8=88\", \"62+30=92\"),\n    @(\"33+54=87\", \"89+3=92\"),\n    @(\"82-58=24\", \"45+45=90\"),\n    @(\"49+11=60\", \"32+35=67\"),\n    @(\"10+86=96\", \"72-8=64\"),\n    @(\"71+19=90\", \"40+14=54\"),\n    @(\"45+32=77\", \"50-43=7\"),\n    @(\"10+53=63\", \"82-32=50\"),\n    @(\"81-16=65\", \"56-8=48\"),\n    @(\"98-41=57\", \"92-40=52\"),\n    @(\"24+1=25\", \"15+52=67\"),\n    @(\"6+85=91\", \"19+11=30\"),\n    @(\"68-11=57\", \"8+34=42\"),\n    @(\"48-19=29\", \"75-63=12\"),\n    @(\"89-46=43\", \"52-8=44\"),\n    @(\"66-32=34\", \"27+32=59\"),\n    @(\"69-1=68\", \"19+66=85\"),\n    @(\"12+87=99\", \"84-73=11\"),\n    @(\"60+27=87\", \"90-53=37\"),\n    @(\"28+71=99\", \"43-9=34\"),\n    @(\"87-74=13\", \"62+8=70\"),\n    @(\"61-6=55\", \"13+83=96\"),\n    @(\"61+3=64\", \"16+4=20\"),\n    @(\"81-45=36\", \"74+12=86\"),\n    @(\"49-35=14\", \"93-73=20\"),\n    @(\"74+5=79\", \"84-82=2\"),\n    @(\"26-7=19\", \"27+28=55\"),\n    @(\"39+46=85\", \"59+24=83\"),\n    @(\"3+44=47\", \"12+1=13\"),\n    @(\"56-0=56\", \"81+2=83\"),\n    @(\"1+13=14\", \"53+35=88\"),\n    @(\"67-60=7\", \"45+30=75\"),\n    @(\"80-4=76\", \"10+67=77\"),\n    @(\"63-41=22\", \"41-12=29\"),\n    @(\"70-32=38\", \"4+15=19\"),\n    @(\"82+11=93\", \"62-14=48\"),\n    @(\"95-38=57\", \"24+10=34\"),\n    @(\"77+10=87\", \"53+25=78\"),\n    @(\"57-0=57\", \"61+22=83\"),\n    @(\"6+6=12\", \"30-10=20\"),\n    @(\"4+76=80\", \"36+32=68\"),\n    @(\"1+87=88\", \"37-11=26\"),\n    @(\"85+11=96\", \"56-36=20\"),\n    @(\"23+5=28\", \"73+13=86\"),\n    @(\"87-75=12\", \"12+54=66\"),\n    @(\"55+14=69\", \"90-34=56\"),\n    @(\"37+17=54\", \"91-75=16\"),\n    @(\"96-5=91\", \"70+18=88\"),\n    @(\"21+45=66\", \"17+4=21\"),\n    @(\"23+28=51\", \"74-48=26\"),\n    @(\"95-75=20\", \"22+50=72\"),\n    @(\"10+44=54\", \"35+29=64\"),\n    @(\"88-47=41\", \"54-35=19\"),\n    @(\"17+32=49\", \"64-33=31\"),\n    @(\"49-33=16\", \"68-50=18\"),\n    @(\"42-21=21\", \"0+11=11\"),\n    @(\"68-41=27\", \"74-10=64\"),\n    @(\"57-54=3\", \"77-70=7\"),\n    @(\"87-34=53\", \"88-72=16\"),\n    @(\"7+60=67\", \"95-14=81\"),\n    @(\"36+12=48\", \"24+14=38\"),\n    @(\"87+1=88\", \"30+57=87\"),\n    @(\"48+35=83\", \"86-67=19\"),\n    @(\"47+49=96\", \"90-0=90\"),\n    @(\"96-92=4\", \"9+42=51\"),\n    @(\"58+37=95\", \"95-11=84\"),\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $rng = $d.Content\n    $found = $rng.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n    if (-not $found) {\n        throw \"Text not found: $old\"\n    }\n}\n"}
